$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to lowercase/underscore convention for consistency
$ws.Range("A1").Value = "indicator_name"
$ws.Range("B1").Value = "actual"
$ws.Range("C1").Value = "actual_lastweek"
$ws.Range("D1").Value = "actual_lastyear"
$ws.Range("E1").Value = "target"
$ws.Range("F1").Value = "perc"
$ws.Range("G1").Value = "perc_week"
$ws.Range("H1").Value = "perc_year"
$ws.Range("I1").Value = "behind_by"
$ws.Range("J1").Value = "text"

# Move the active selection to J1
$ws.Range("J1").Select()
